$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing existing rows 104-160 down to 105-161
$ws.Rows.Item(104).Insert()

# Fill the new row 104 with the weekly price record
$ws.Cells.Item(104, 1).Value = 4
$ws.Cells.Item(104, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value = "Los Lagos"
$ws.Cells.Item(104, 4).Value = 45176
$ws.Cells.Item(104, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 5).Value = 10
$ws.Cells.Item(104, 6).Value = 100112026
$ws.Cells.Item(104, 7).Value = "Haba"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 40
$ws.Cells.Item(104, 11).Value = 16000
$ws.Cells.Item(104, 12).Value = 16000
$ws.Cells.Item(104, 13).Value = 16000
$ws.Cells.Item(104, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(104, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(104, 16).Value = 640
$ws.Cells.Item(104, 17).Value = 25
$ws.Cells.Item(104, 18).Value = "Hortaliza"
